# Update the EPEX Spot price workbook:
#  - "Prix Spot" sheet: add a new day column AS ("28-jul") with hourly prices
#  - "Gaz" sheet: append two new daily rows (2025-07-26, 2025-07-27)
#  - "CO2" sheet: append two new daily rows (2025-07-26, 2025-07-27)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column AS = "28-jul"
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (AR1) onto the new header
# cell (AS1) so it keeps the bold / bordered / centered header style, then
# set its text.
$prevHeader = $wsSpot.Range("AR1")
$prevHeader.Copy()
$newHeader = $wsSpot.Range("AS1")
$newHeader.PasteSpecial(-4122)
$newHeader.Value = "28-jul"

$spotValues = @{
    2  = "62.29"
    3  = "55.75"
    4  = "48.95"
    5  = "19.96"
    6  = "38"
    7  = "45.07"
    8  = "50"
    9  = "56.23"
    10 = "63.71"
    11 = "30"
    12 = "20.59"
    13 = "22.49"
    14 = "21.87"
    15 = "20.22"
    16 = "11.17"
    17 = "10.87"
    18 = "14.35"
    19 = "22.97"
    20 = "44.3"
    21 = "60.11"
    22 = "60"
    23 = "74.98"
    24 = "88.54000000000001"
    25 = "63.25"
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item([int]$row, 45).Value = [double]$spotValues[$row]
}

# ---------------------------------------------------------------------------
# Helper: append a (date-as-text, value) row to a two-column sheet while
# keeping the date stored as plain text (not auto-converted to a date
# serial number) and without leaving behind any extra cell style.
# ---------------------------------------------------------------------------
function Add-DailyRow($ws, $rowIndex, $dateText, $value) {
    $dateCell = $ws.Cells.Item($rowIndex, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.Style = "Normal"
    $ws.Cells.Item($rowIndex, 2).Value = $value
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append rows 42 and 43
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
Add-DailyRow $wsGaz 42 "2025-07-26" 31.85
Add-DailyRow $wsGaz 43 "2025-07-27" 31.85

# ---------------------------------------------------------------------------
# Sheet "CO2": append rows 42 and 43
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
Add-DailyRow $wsCo2 42 "2025-07-26" 70.7
Add-DailyRow $wsCo2 43 "2025-07-27" 70.7
